$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B2's text (value itself is unchanged: "Tuple[int] (optional)") but
# re-assigning .Value resets the cell's quotePrefix style (s="1"), so stash
# the original formatting in a scratch cell, write the value, then restore it.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B2").Value = "Tuple[int] (optional)"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("Z1").Clear() | Out-Null

# --- Correct the typo in C5's text: "False" -> "False."
$ws.Range("C5").Value = " If True, it only bads to right and bottom. Defaults to False."

# --- Move the active selection to C6 (matches the saved cursor position)
$ws.Range("C6").Select()
